$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Values that read as plain numbers (e.g. "242.62") get a leading apostrophe so
# Excel stores them as literal text, same as the original scraped/inline strings,
# instead of silently coercing them into numeric cells.

$ws.Range("D2").Value = "29.328.89"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "1.877.75"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'0.7259"
$ws.Range("E5").Value = "  +2.47%  "

$ws.Range("D6").Value = "'242.62"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "'0.08000"
$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("D9").Value = "'0.3160"
$ws.Range("E9").Value = "  +2.22%  "

$ws.Range("D10").Value = "'24.97"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").Value = "'0.08223"

$ws.Range("D12").Value = "1.887.76"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").Value = "'94.77"
$ws.Range("E13").Value = "  +4.14%  "

$ws.Range("D14").Value = "'5.227"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "'0.7122"
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").Value = "'6.389"
$ws.Range("E16").Value = "  +5.25%  "

$ws.Range("D17").Value = "'0.000008486"
$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").Value = "29.333.30"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "'243.53"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").Value = "2.138.38"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").Value = "'7.775"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").Value = "'0.1607"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D26").Value = "'162.65"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").Value = "'9.028"
$ws.Range("E27").Value = "  +0.27%  "

$ws.Range("D28").Value = "'18.51"
$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("D29").Value = "'1.503"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'4.407"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "'4.308"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").Value = "'1.190"
$ws.Range("E32").Value = "  -7.69%  "

$ws.Range("D33").Value = "'0.05356"
$ws.Range("E33").Value = "  +0.40%  "

$ws.Range("D34").Value = "'1.941"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("D35").Value = "'0.7588"
$ws.Range("E35").Value = "  +1.91%  "

$ws.Range("D36").Value = "'1.176"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").Value = "'2.696"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").Value = "'0.01879"
$ws.Range("E38").Value = "  +0.58%  "

$ws.Range("D39").Value = "1.278.42"
$ws.Range("E39").Value = "  +3.86%  "

$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("D41").Value = "'6.431"
$ws.Range("E41").Value = "  -2.00%  "

$ws.Range("D42").Value = "'113.22"
$ws.Range("E42").Value = "  +3.49%  "

$ws.Range("D43").Value = "'0.9088"
$ws.Range("E43").Value = "  +2.57%  "

$ws.Range("D44").Value = "'74.33"
$ws.Range("E44").Value = "  +2.70%  "

$ws.Range("D45").Value = "'0.00000000133"
$ws.Range("E45").Value = "  +8.67%  "

$ws.Range("D47").Value = "2.034.06"
$ws.Range("E47").Value = "  +0.97%  "

$ws.Range("D48").Value = "'0.5229"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "'1.795"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").Value = "'9.492"
$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("D51").Value = "'0.4346"
$ws.Range("E51").Value = "  +0.86%  "
